$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 83334536
$ws.Range("I11").Value = 83334536
$ws.Range("K11").Value = 83334536
$ws.Range("M11").Value = -83334396

$ws.Range("H19").Value = 1770
$ws.Range("J19").Value = 1729.6
$ws.Range("L19").Value = 1729.6
$ws.Range("N19").Value = -2079.6

$ws.Range("H98").Value = 1041.1321
$ws.Range("I98").Value = 926.53845
$ws.Range("K98").Value = 926.53845
$ws.Range("M98").Value = 571.46155

$ws.Range("H122").Value = 1041.1321
$ws.Range("I122").Value = 926.53845
$ws.Range("K122").Value = 2779.61535
$ws.Range("M122").Value = -329.61535

$ws.Range("H137").Value = 67152.38
$ws.Range("I137").Value = 93843.664
$ws.Range("J137").Value = 7097
$ws.Range("K137").Value = 281530.992
$ws.Range("L137").Value = 21291
$ws.Range("M137").Value = -278980.992
$ws.Range("N137").Value = -26391

$ws.Range("H138").Value = 2954.6
$ws.Range("I138").Value = 1651.7778
$ws.Range("J138").Value = 3823.1482
$ws.Range("K138").Value = 4955.3334
$ws.Range("L138").Value = 11469.4446
$ws.Range("M138").Value = 184.6665999999996
$ws.Range("N138").Value = -21749.4446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3697.158
$ws.Range("I32").Value = 2881.7058
$ws.Range("J32").Value = 10628.5
$ws.Range("K32").Value = 2881.7058
$ws.Range("L32").Value = 10628.5
$ws.Range("M32").Value = -2594.7058
$ws.Range("N32").Value = -11202.5

$ws.Range("H74").Value = 1712.75
$ws.Range("I74").Value = 1645.0526
$ws.Range("J74").Value = 2999
$ws.Range("K74").Value = 1645.0526
$ws.Range("L74").Value = 2999
$ws.Range("M74").Value = -771.0526
$ws.Range("N74").Value = -4747

$ws.Range("H77").Value = 1712.75
$ws.Range("I77").Value = 1645.0526
$ws.Range("J77").Value = 2999
$ws.Range("K77").Value = 8225.262999999999
$ws.Range("L77").Value = 14995
$ws.Range("M77").Value = -3857.262999999999
$ws.Range("N77").Value = -23731

$ws.Range("H102").Value = 2306
$ws.Range("I102").Value = 2065.3333
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 2065.3333
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -443.3332999999998
$ws.Range("N102").Value = -6994

$ws.Range("H103").Value = 65686.57000000001
$ws.Range("J103").Value = 65686.57000000001
$ws.Range("L103").Value = 65686.57000000001
$ws.Range("N103").Value = -68030.57000000001

$ws.Range("H132").Value = 325484.1
$ws.Range("I132").Value = 359855.97
$ws.Range("J132").Value = 4680
$ws.Range("K132").Value = 1079567.91
$ws.Range("L132").Value = 14040
$ws.Range("M132").Value = -1077037.91
$ws.Range("N132").Value = -19100

$ws.Range("H133").Value = 161902.5
$ws.Range("J133").Value = 161902.5
$ws.Range("L133").Value = 161902.5
$ws.Range("N133").Value = -166962.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 105102
$ws.Range("J60").Value = 105102
$ws.Range("L60").Value = 105102
$ws.Range("N60").Value = -106300

$ws.Range("H105").Value = 3595.8462
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.642857
$ws.Range("J7").Value = 12
$ws.Range("L7").Value = 12
$ws.Range("N7").Value = -238

$ws.Range("H58").Value = 3318.4707
$ws.Range("I58").Value = 3021.4783
$ws.Range("J58").Value = 3939.4546
$ws.Range("K58").Value = 3021.4783
$ws.Range("L58").Value = 3939.4546
$ws.Range("M58").Value = -2818.4783
$ws.Range("N58").Value = -4345.4546

$ws.Range("H132").Value = 5329.357
$ws.Range("I132").Value = 5161.3
$ws.Range("K132").Value = 15483.9
$ws.Range("M132").Value = -12953.9

$ws.Range("H134").Value = 4916.1665
$ws.Range("I134").Value = 3999.4285
$ws.Range("K134").Value = 11998.2855
$ws.Range("M134").Value = -9463.2855

$ws.Range("H136").Value = 3318.4707
$ws.Range("I136").Value = 3021.4783
$ws.Range("J136").Value = 3939.4546
$ws.Range("K136").Value = 9064.4349
$ws.Range("L136").Value = 11818.3638
$ws.Range("M136").Value = -6514.4349
$ws.Range("N136").Value = -16918.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 569
$ws.Range("I5").Value = 533.0833
$ws.Range("K5").Value = 1599.2499
$ws.Range("M5").Value = -1487.2499

$ws.Range("H11").Value = 13097.25
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H14").Value = 17371.428
$ws.Range("I14").Value = 17371.428
$ws.Range("K14").Value = 52114.284
$ws.Range("M14").Value = -51941.284

$ws.Range("H22").Value = 3008.75
$ws.Range("I22").Value = 3008.75
$ws.Range("K22").Value = 9026.25
$ws.Range("M22").Value = -8857.25

$ws.Range("H26").Value = 1066.6666
$ws.Range("I26").Value = 1000
$ws.Range("J26").Value = 1200
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3600
$ws.Range("M26").Value = -2712
$ws.Range("N26").Value = -4176

$ws.Range("H27").Value = 3008.75
$ws.Range("I27").Value = 3008.75
$ws.Range("K27").Value = 9026.25
$ws.Range("M27").Value = -8924.25

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H87").Value = 175
$ws.Range("I87").Value = 175
$ws.Range("K87").Value = 525
$ws.Range("M87").Value = 723

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H90").Value = 175
$ws.Range("I90").Value = 175
$ws.Range("K90").Value = 1575
$ws.Range("M90").Value = 4665

$ws.Range("H135").Value = 569
$ws.Range("I135").Value = 533.0833
$ws.Range("K135").Value = 4797.7497
$ws.Range("M135").Value = -2262.7497

$ws.Range("H139").Value = 2892.25
$ws.Range("I139").Value = 2674.7222
$ws.Range("K139").Value = 8024.1666
$ws.Range("M139").Value = -2884.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4394.706
$ws.Range("I126").Value = 3631.9
$ws.Range("K126").Value = 10895.7
$ws.Range("M126").Value = -8425.700000000001

$ws.Range("H132").Value = 3268.762
$ws.Range("I132").Value = 3268.762
$ws.Range("K132").Value = 9806.286
$ws.Range("M132").Value = -7276.286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1912
$ws.Range("J16").Value = 2050
$ws.Range("L16").Value = 2050
$ws.Range("N16").Value = -2390

$ws.Range("H55").Value = 294.41666
$ws.Range("I55").Value = 334
$ws.Range("J55").Value = 269.22726
$ws.Range("K55").Value = 334
$ws.Range("L55").Value = 269.22726
$ws.Range("M55").Value = -161
$ws.Range("N55").Value = -615.22726

$ws.Range("H93").Value = 111112570
$ws.Range("J93").Value = 1976
$ws.Range("L93").Value = 1976
$ws.Range("N93").Value = -4472

$ws.Range("H122").Value = 90668
$ws.Range("I122").Value = 36002
$ws.Range("K122").Value = 108006
$ws.Range("M122").Value = -105556

$ws.Range("H132").Value = 912417.2
$ws.Range("I132").Value = 1114287.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 3342862.8
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3340332.8
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 4198.357
$ws.Range("I136").Value = 3940.2
$ws.Range("J136").Value = 4843.75
$ws.Range("K136").Value = 11820.6
$ws.Range("L136").Value = 14531.25
$ws.Range("M136").Value = -9270.599999999999
$ws.Range("N136").Value = -19631.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 16431.182
$ws.Range("I96").Value = 9332.5
$ws.Range("K96").Value = 9332.5
$ws.Range("M96").Value = -7959.5

$ws.Range("H132").Value = 25529.283
$ws.Range("I132").Value = 32121.97
$ws.Range("J132").Value = 4552.5454
$ws.Range("K132").Value = 96365.91
$ws.Range("L132").Value = 13657.6362
$ws.Range("M132").Value = -93835.91
$ws.Range("N132").Value = -18717.6362

$ws.Range("H136").Value = 20302.055
$ws.Range("I136").Value = 1342.9756
$ws.Range("J136").Value = 80096.08
$ws.Range("K136").Value = 4028.9268
$ws.Range("L136").Value = 240288.24
$ws.Range("M136").Value = -1478.9268
$ws.Range("N136").Value = -245388.24
